$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 108, shifting existing rows 108-118 down to 109-119.
$ws.Rows.Item(108).Insert()

# Populate the newly inserted row 108 with its data (weekly update row).
$ws.Cells.Item(108, 1).Value = 7
$ws.Cells.Item(108, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(108, 3).Value = "Ñuble"
$ws.Cells.Item(108, 4).Value = 44449
$ws.Cells.Item(108, 5).Value = 16
$ws.Cells.Item(108, 6).Value = 100112006
$ws.Cells.Item(108, 7).Value = "Repollo"
$ws.Cells.Item(108, 8).Value = "Crespo record"
$ws.Cells.Item(108, 9).Value = "Primera"
$ws.Cells.Item(108, 10).Value = 300
$ws.Cells.Item(108, 11).Value = 600
$ws.Cells.Item(108, 12).Value = 650
$ws.Cells.Item(108, 13).Value = 625
$ws.Cells.Item(108, 14).Value = "$/unidad"
$ws.Cells.Item(108, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(108, 16).Value = 625
$ws.Cells.Item(108, 17).Value = 1
$ws.Cells.Item(108, 18).Value = "Hortaliza"
